# Re-ran analyses with slight plotting updates.
# This updates the "meta analyses" sheet: the old "Verification attempt 3"
# (Hunter & Schmidt implementation 2 / Field & Gillett Basic meta-analysis)
# row is dropped, and the two remaining verification rows are renumbered
# down by one ("attempt 4" -> "attempt 3", "attempt 5" -> "attempt 4").

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("power analyses")
$ws2 = $wb.Worksheets.Item("meta analyses")

# ---------------------------------------------------------------------
# "meta analyses" sheet: remove the old row 6 and rename the two
# verification rows that shift up into rows 6 and 7.
# ---------------------------------------------------------------------
$ws2.Rows.Item(6).Delete() | Out-Null

$ws2.Range("A6").Value = "Verification attempt 3: Hunter & Schmidt method (implementation 3)"
$ws2.Range("A7").Value = "Verification attempt 4: Mix of Hunter & Schmidt and Hedges' methods"

# ---------------------------------------------------------------------
# Cosmetic cleanup of redundant direct-formatting entries left over on
# "power analyses" after the re-run (no visible change: these fills were
# all "no fill"/default, Excel just collapses the duplicate style refs).
# ---------------------------------------------------------------------
$ws1.Range("D3:D6").NumberFormat = "0.00"
$ws1.Range("G3:G6").NumberFormat = "0.00"
$ws1.Range("G7:G10").NumberFormat = "0.00"
$ws1.Range("H3:H10").ClearFormats()
$ws1.Range("D7:D10").ClearFormats()

# ---------------------------------------------------------------------
# Restore selections on both sheets to match where the author ended up.
# ---------------------------------------------------------------------
$ws1.Activate() | Out-Null
$ws1.Range("A7").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("B15").Select() | Out-Null
